$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "(population census results)" sub-heading row (row 2).
# The rows below shift up: old row3(blank)->2, row4("(sq km)")->3,
# row5(years)->4, row6(area figures)->5.
$ws.Rows(2).Delete()

# Remove the 1989 and 2002 columns, keeping only the 2014 column
# (which becomes column B).
$ws.Columns(2).Delete()
$ws.Columns(2).Delete()

# The workbook/sheet used to just be called "1" - give it the real
# municipality name matching the file name.
$ws.Name = "კასპი"

$ws.Range("A2").Select() | Out-Null
